# Applied json validation to addProduct API request body
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New requirement note in row 4 (To Learn column)
$ws.Range("E4").Value = "JSON Schema validator"

# Row 6 becomes a fresh/blank requirement row - the requirement that used
# to live here moves down into row 7, so wipe row 6's A:C content+format.
$ws.Range("A6:C6").Clear()
$ws.Rows.Item(6).RowHeight = 19.5
$ws.Rows.Item(7).RowHeight = 32.25

# Row 7 now holds what used to be row 6's requirement (already DONE-styled)
$ws.Range("A7").Value = "Clean code over all"
$ws.Range("B7").Value = "productController.js & product.model"

# Row 8 now holds what used to be row 7's requirement, newly marked DONE
$ws.Range("A8").Value = "Draw block diagram of req, res cycle as per my understandng"
$ws.Range("C2").Copy()
$ws.Range("C8").PasteSpecial(-4122)
$ws.Range("C8").Value = "DONE"

# Row 9 now holds what used to be row 8's requirement; its old B-column
# note is cleared (format reset to the plain label style) and it is
# marked DONE
$ws.Range("A9").Value = "Checkout NoSQL booster"
$ws.Range("B9").ClearContents()
$ws.Range("B8").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("C9").Value = "DONE"

# Row 10 now holds what used to be row 9's requirement and its note,
# newly marked DONE
$ws.Range("A10").Value = "Add pagination to getProducts API of productController.js"
$ws.Range("B2").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("B10").Value = "productController.js"
$ws.Range("C2").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$ws.Range("C10").Value = "DONE"

# New requirement added in row 12
$ws.Range("A12").Value = "vs code extension : code spell checker"

# Update the saved selection and page orientation
$null = $ws.Range("A1:F12").Select()
$ws.PageSetup.Orientation = 1
